$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("pages_with_footer")
$ws2 = $wb.Worksheets.Item("pages_with_footer_spanish")

# Delete the Spanish worksheet - its rows are being folded into sheet1
$excel.DisplayAlerts = $false
$ws2.Delete()
$excel.DisplayAlerts = $true

# Rebuild sheet1 data: combine English + Spanish providers into a single table
# and add a couple of "Blog" language/footer test rows.
$ws1.Range("A1").Value = "path"
$ws1.Range("B1").Value = "type"
$ws1.Range("C1").Value = "language"

$ws1.Range("A2").Value = "about-cancer/coping/feelings/relaxation/vitamin-d-supplement-cancer-prevention"
$ws1.Range("B2").Value = "Blog Page"
$ws1.Range("C2").Value = "English"

$ws1.Range("A3").Value = "/espanol/node/46"
$ws1.Range("B3").Value = "Press Release"
$ws1.Range("C3").Value = "Spanish"

$ws1.Range("A4").Value = "about-cancer/coping/feelings"
$ws1.Range("B4").Value = "Article"
$ws1.Range("C4").Value = "English"

$ws1.Range("A5").Value = "/node/36"
$ws1.Range("B5").Value = "Cancer Center"
$ws1.Range("C5").Value = "English"

$ws1.Range("A6").Value = "types/breast/patient/breast-treatment-pdq"
$ws1.Range("B6").Value = "PDQ Cancer Information Summary"
$ws1.Range("C6").Value = "English"

$ws1.Range("A7").Value = "/news-events/press-releases/2018/oropharyngeal-hpv-cisplatin"
$ws1.Range("B7").Value = "Article"
$ws1.Range("C7").Value = "English"

$ws1.Range("A8").Value = "espanol/cancer/sobrellevar/sentimientos/hoja-informativa-estres"
$ws1.Range("B8").Value = "Article"
$ws1.Range("C8").Value = "Spanish"

$ws1.Range("A9").Value = "about-cancer/coping/feelings/relaxation/hpv-vaccine-presidents-cancer-panel-improving-uptake"
$ws1.Range("B9").Value = "Blog"
$ws1.Range("C9").Value = "English"

$ws1.Range("A10").Value = "about-cancer/coping/feelings"
$ws1.Range("B10").Value = "Article"
$ws1.Range("C10").Value = "English"

# Wrap text on the final row's path cell (long URL path column)
$ws1.Range("A10").WrapText = $true

# Update selection/view state
$ws1.Range("A13").Select()
$excel.ActiveWindow.WindowState = $excel.ActiveWindow.WindowState
